# Apply the table-style change recorded on slide 6: the data table there
# switches from the deck's default table style to the built-in
# "Medium Style 2 - Accent 1" style ({C8549F4C-F7FD-48B3-98FA-6FC13F9D7C89}).

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(6)

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{C8549F4C-F7FD-48B3-98FA-6FC13F9D7C89}")
    }
}
